$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.17532746519896136
$ws.Range("A2").Value = -0.0099999995554611587
$ws.Range("A3").Value = -0.0089999995563676549
$ws.Range("A4").Value = 0.20999999883007092
$ws.Range("A5").Value = -0.0059999995664359318
$ws.Range("A6").Value = -0.0059999995485391366
$ws.Range("A7").Value = -0.019999999468049978
$ws.Range("A8").Value = -0.048139182020049631
$ws.Range("A9").Value = -0.0059999995407382656
$ws.Range("A10").Value = -0.0059999995396182726
$ws.Range("A11").Value = -0.0044999995479706456
$ws.Range("A12").Value = -0.0059999995395481065
$ws.Range("A13").Value = -0.0059999995403128281
$ws.Range("A14").Value = -0.011999999507151138
$ws.Range("A15").Value = 0.045457563862300709
$ws.Range("A16").Value = -0.0059999995408310802
$ws.Range("A17").Value = -0.0059999995389210525
$ws.Range("A18").Value = -0.0089999995215999107
$ws.Range("A19").Value = -0.047188376518308761
$ws.Range("A20").Value = -0.0089999995562397572
$ws.Range("A21").Value = -0.0089999995556766521
$ws.Range("A22").Value = -0.0089999995552707546
$ws.Range("A23").Value = -0.0089999995525884557
$ws.Range("A24").Value = -0.041999999359275009
$ws.Range("A25").Value = -0.041999999355874174
$ws.Range("A26").Value = -0.0059999995465958023
$ws.Range("A27").Value = -0.005999999544219925
$ws.Range("A28").Value = -0.00599999953410979
$ws.Range("A29").Value = -0.011999999493504276
$ws.Range("A30").Value = -0.019999999445131422
$ws.Range("A31").Value = -0.014999999468427561
$ws.Range("A32").Value = -0.0209999994337986
$ws.Range("A33").Value = -0.0059999995180826105
